$p = $ppt.ActivePresentation

# 1. Merge the "Chapitre " + "15" runs on slide 1 into a single run "Chapitre 15"
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$sub = $tr.Characters(1, 11)
$sub.Text = "Chapitre 15"

# 2. Delete the last slide (slide 17, "URL Rewriting")
$p.Slides.Item(17).Delete()
